# Delete the "SAMSUNG_GalaxyM02_Android_11.0.0_51323" device column (column F)
# from the "DeviceList" sheet, shifting subsequent columns left by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

$ws.Columns("F").Delete()

# The conditional formatting ranges used to span B2:I2 (including the now
# removed column); re-point them at the shrunk B2:H2 range.
$newRange = $ws.Range("B2:H2")
$fcs = $newRange.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($newRange)
}

# Update selection on this sheet to match the post-edit state.
$ws.Activate()
$ws.Range("E15").Select()
